# Generate Report for Handback
#
# Two files ("53eb5cf7-...md" itself, and the file that depends on it,
# "95e71b23-...md") have now been handed back for both the zh-cn and
# de-de locales. Update the Overview + per-locale status sheets to
# reflect the handback: flip Status from "Ready for handoff" to
# "Handed back: in sync with en-US", and record the Latest Target
# File / Latest Handback File / Latest Handback DateTime for each
# locale sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: rows for 53eb5cf7-....md (row 3) and
# 95e71b23-....md (row 4) both move from "Ready for handoff" to
# "Handed back: in sync with en-US" for both zh-cn and de-de columns.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value2 = "Handed back: in sync with en-US"
$overview.Range("C3").Value2 = "Handed back: in sync with en-US"
$overview.Range("B4").Value2 = "Handed back: in sync with en-US"
$overview.Range("C4").Value2 = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B3").Value2 = "Handed back: in sync with en-US"
$zh.Range("B4").Value2 = "Handed back: in sync with en-US"

$zh.Range("E3").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md"
$zh.Range("F3").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.zh-cn.xlf"
$zh.Range("G3").Value2 = "2016-03-07 02:24:30"

$zh.Range("E4").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md"
$zh.Range("F4").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.zh-cn.xlf"
$zh.Range("G4").Value2 = "2016-03-07 02:24:30"

$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.zh-cn.xlf", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md")
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.zh-cn.xlf", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.zh-cn.xlf")

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B3").Value2 = "Handed back: in sync with en-US"
$de.Range("B4").Value2 = "Handed back: in sync with en-US"

$de.Range("E3").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md"
$de.Range("F3").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.de-de.xlf"
$de.Range("G3").Value2 = "2016-03-07 02:24:48"

$de.Range("E4").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md"
$de.Range("F4").Value2 = "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.de-de.xlf"
$de.Range("G4").Value2 = "2016-03-07 02:24:48"

$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.de-de.xlf", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.md")
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.de-de.xlf", "", "", "53eb5cf7-86ad-4a68-9abf-c9f7b4b7c933.40a5f0171868f5649ac3e1cc03d66877e85c1eb3.de-de.xlf")

Write-Host "Handback report generated."
